$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.801.76"
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("D3").Value = "3.056.58"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'518.19"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").Value = "'142.29"
$ws.Range("E6").Value = "  +3.18%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.434"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("D9").Value = "'7.26"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("D12").Value = "3.587.96"
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").Value = "'25.92"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "57.831.98"
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("D17").Value = "3.061.36"
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").Value = "'12.82"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Value = "'8.12"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'330.46"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "'0.498"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "'65.71"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  -4.02%  "
$ws.Range("D28").Value = "'6.39"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Value = "'7.22"
$ws.Range("E29").Value = "  +4.46%  "
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("E31").Value = "  +2.75%  "
$ws.Range("D32").Value = "'20.63"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").Value = "'154.44"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").Value = "'4.51"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").Value = "'27.12"
$ws.Range("E35").Value = "  +3.95%  "
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("D39").Value = "3.102.31"
$ws.Range("E39").Value = "  +2.41%  "
$ws.Range("D40").Value = "'3.91"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("D41").Value = "'36.54"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'0.654"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Value = "2.266.69"
$ws.Range("E44").Value = "  +3.53%  "
$ws.Range("E45").Value = "  +9.94%  "
$ws.Range("D46").Value = "'20.73"
$ws.Range("E46").Value = "  +6.34%  "
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").Value = "'5.87"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").Value = "'0.927"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "'0.729"
$ws.Range("E50").Value = "  +8.20%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'258.83"
$ws.Range("E51").Value = "  +11.93%  "
